# AllianceShrine.xlsx — "dev of alliance village"
#
# The shrine-stage troop composition for STR_troops (column I) is unified
# across every stage row: all rows now reference the same troop string,
# which also introduces a red-dragon garrison at the first slot.
#
# This naturally prunes the old, now-unreferenced per-stage troop strings
# from the shared string table (they were only ever referenced from this
# column), matching the shared-strings shrink seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTroops = "redDragon_1_5,swordsman_1_10,sentinel_1_10,ranger_1_10&swordsman_1_10,sentinel_1_10,ranger_1_10"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 9).Value = $newTroops
}

# Reflect the author's final selection/view state on the sheet.
[void]$ws.Range("I25").Select()
